$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Days in Shop
$ws.Range("I2").Value = 35

# Row 3: Production Date and Days in Shop
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2024-07-28"
$ws.Range("I3").Value = 8

# Row 4: Production Date and Days in Shop
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2024-08-02"
$ws.Range("I4").Value = 3

# Row 5: full row of data corrected
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2024-03-14"
$ws.Range("B5").Value = "Walter"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "92020483"
$ws.Range("D5").Value = "Aomdonm"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2024-05-28"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "$3920"
$ws.Range("G5").Value = "Ansoansoansaonason"
$ws.Range("H5").Value = "J99999"
$ws.Range("I5").Value = 94
